# Remove the gray (bg1, lumMod 50%) outline from the themed code-block
# placeholders on the "Title and Content" and "Two Content" slide layouts,
# replacing the <a:solidFill> line color with <a:noFill/> (line stays the
# same weight, just becomes invisible/unfilled).

$p = $ppt.ActivePresentation

# Layout 2: "Title and Content"  -> shapes "Content Placeholder 2", "Title 6"
# Layout 4: "Two Content"        -> shapes "Title 1", "Content Placeholder 2",
#                                    "Content Placeholder 3"
$targets = @{
    "Title and Content" = @("Content Placeholder 2", "Title 6")
    "Two Content"       = @("Title 1", "Content Placeholder 2", "Content Placeholder 3")
}

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $lay = $p.SlideMaster.CustomLayouts.Item($li)
    if ($targets.ContainsKey($lay.Name)) {
        $names = $targets[$lay.Name]
        for ($si = 1; $si -le $lay.Shapes.Count; $si++) {
            $sh = $lay.Shapes.Item($si)
            if ($names -contains $sh.Name) {
                $sh.Line.Visible = 0
            }
        }
    }
}
